# Apply updated crypto price/volume figures to worksheet cells,
# as scraped on Thu Aug  3 10:12:33 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.124.62'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '1.833.02'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.88'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6644'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9997'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2949'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07332'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -4.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.73'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07670'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').Value = '1.843.39'
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('E13').Value = '  -2.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6743'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '86.14'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -5.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.106'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.70%  '
$ws.Range('D17').Value = '29.117.40'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008231'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '227.91'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -4.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.49'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9989'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.282'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.15%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '160.78'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('E25').Value = '  -4.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.662'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.97'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.499'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.236'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.095'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.197'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05336'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.864'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7442'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.128'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.70%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').Value = '1.312.25'
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01801'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.709'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9263'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.019'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9978'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '103.33'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.76%  '
$ws.Range('D44').Value = '1.980.05'
$ws.Range('E44').Value = '  -1.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5167'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('E46').Value = '  -3.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.07668'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +14.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.759'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '63.29'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.264'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05924'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.08%  '
